$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values in column D must stay as text (inline strings),
# so format the cells as Text before assigning to avoid Excel auto-converting
# them to numbers (which would also normalize/strip formatting like trailing zeros).
$dCells = @("D2","D4","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D27","D40","D41","D42","D43","D44","D45","D48")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.34"
$ws.Range("D4").Value = "5.417"
$ws.Range("D6").Value = "3.454"
$ws.Range("D7").Value = "6.537"
$ws.Range("D9").Value = "0.9123"
$ws.Range("D10").Value = "0.1408"
$ws.Range("D11").Value = "0.07483"
$ws.Range("D12").Value = "0.03299"
$ws.Range("D13").Value = "0.03057"
$ws.Range("D14").Value = "0.09347"
$ws.Range("D15").Value = "3.857"
$ws.Range("D16").Value = "0.001562"
$ws.Range("D17").Value = "0.04669"
$ws.Range("D18").Value = "0.0005942"
$ws.Range("D19").Value = "0.006121"
$ws.Range("D20").Value = "0.004988"
$ws.Range("D21").Value = "0.0009807"
$ws.Range("D22").Value = "0.0001101"
$ws.Range("D23").Value = "3.605"
$ws.Range("D24").Value = "2.135"
$ws.Range("D27").Value = "0.0002901"
$ws.Range("D40").Value = "0.03960"
$ws.Range("D41").Value = "0.006216"
$ws.Range("D42").Value = "0.1077"
$ws.Range("D43").Value = "0.002621"
$ws.Range("D44").Value = "0.008698"
$ws.Range("D45").Value = "0.00005230"
$ws.Range("D48").Value = "0.8983"

# Plain text fields (coin names, links, volume labels) can be set directly.
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
